$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.813.86'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.145.58'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.77'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.55'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.145.67'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.11'
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.496'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.91'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.658.47'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.919.77'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.144.70'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '501.26'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.74'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.710'
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.11'
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.77'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.83'
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.79'
$ws.Range("E30").Value = '  +5.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.42'
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.15'
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  -2.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.68'
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0894'
$ws.Range("E37").Value = '  +4.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '472.20'
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0414'
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.61'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.013.51'
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("E43").Value = '  -3.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.280'
$ws.Range("E44").Value = '  -3.59%  '
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.07'
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("E50").Value = '  -3.44%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.03'
$ws.Range("E51").Value = '  +4.66%  '
